$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6043246984481812
$ws.Range("B1").Value = 4.108129978179932
$ws.Range("C1").Value = 6.178332328796387
$ws.Range("D1").Value = 1.505155205726624
$ws.Range("E1").Value = 0.8473306894302368
